# Advertising verbose spreadsheet update (2026-01-14 generated outputs)
#
# The "documents" block (Application -> documents -> file, rows 12-15) gains a
# new leading field "uploaded-date" (row 12), pushing the existing file fields
# down by one row. Everything from the old row 12 through the old row 23 shifts
# down by one row, and the field that used to be the last row of that block
# (document-reference / name, old row 23) is removed entirely so the table
# keeps the same overall number of rows (the later "Advertisement location"
# etc. rows stay put). Separately, row 79's datatype changes from "string" to
# "enum".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row above the old row 12 ("file"/base64-content row). This
#    shifts rows 12-96 down to 13-97 and correctly grows/shifts the existing
#    merged cells (A2:A18 -> A2:A19, A19:A24 -> A20:A25, etc.).
$ws.Rows.Item(12).Insert()

# 2) Populate the new row 12 with the "uploaded-date" field. Columns A/B stay
#    blank (this is a continuation row of the "Application" / "documents"
#    group, same as the old row 12 it now sits above of). Columns C-G repeat
#    the same breadcrumb as the surrounding "documents" rows.
$ws.Range("C12").Value2 = "The details of the application payload to be submitted"
$ws.Range("D12").Value2 = "application"
$ws.Range("E12").Value2 = "Application"
$ws.Range("F12").Value2 = "documents"
$ws.Range("G12").Value2 = "Documents[]"
$ws.Range("H12").Value2 = "uploaded-date"
$ws.Range("I12").Value2 = "Uploaded date"
$ws.Range("L12").Value2 = "The date the document was uploaded to the application"
$ws.Range("M12").Value2 = "date"
$ws.Range("N12").Value2 = "MUST"

# 3) The row that used to be row 23 ("document-reference" / "name" - "A name
#    for the document. For example, The Site Plan") has now shifted to row 24
#    and is dropped entirely from the table, so the rows below it ("is-advert-
#    overhanging" etc.) snap back up to their original row numbers and the
#    "Advertisement location" merge shrinks back from A20:A25 to A20:A24.
$ws.Rows.Item(24).Delete()

# 4) Unrelated datatype correction further down the sheet: "Illumination
#    method" changes from a free-text "string" to an "enum".
$ws.Range("M79").Value2 = "enum"
